# Auto-generated edit script: apply updated Leve profit figures per sheet.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 25102314
$ws.Range("I132").Value = 28687194
$ws.Range("K132").Value = 86061582
$ws.Range("M132").Value = -86059052
$ws.Range("H135").Value = 407.46155
$ws.Range("I135").Value = 358.08334
$ws.Range("K135").Value = 3222.75006
$ws.Range("M135").Value = -687.7500600000003
$ws.Range("H138").Value = 2674.2144
$ws.Range("I138").Value = 1307.6364
$ws.Range("J138").Value = 3558.4707
$ws.Range("K138").Value = 3922.9092
$ws.Range("L138").Value = 10675.4121
$ws.Range("M138").Value = 1217.0908
$ws.Range("N138").Value = -20955.4121
$ws.Range("H141").Value = 3815.4167
$ws.Range("I141").Value = 3611.875
$ws.Range("J141").Value = 4222.5
$ws.Range("K141").Value = 10835.625
$ws.Range("L141").Value = 12667.5
$ws.Range("M141").Value = -5655.625
$ws.Range("N141").Value = -23027.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 644.1539
$ws.Range("I2").Value = 579.5454999999999
$ws.Range("J2").Value = 999.5
$ws.Range("K2").Value = 579.5454999999999
$ws.Range("L2").Value = 999.5
$ws.Range("M2").Value = -466.5454999999999
$ws.Range("N2").Value = -1225.5
$ws.Range("H32").Value = 6274.4287
$ws.Range("I32").Value = 4324.512
$ws.Range("J32").Value = 9908.362999999999
$ws.Range("K32").Value = 4324.512
$ws.Range("L32").Value = 9908.362999999999
$ws.Range("M32").Value = -4037.512
$ws.Range("N32").Value = -10482.363
$ws.Range("H37").Value = 30723.834
$ws.Range("I37").Value = 29800
$ws.Range("J37").Value = 30908.6
$ws.Range("K37").Value = 29800
$ws.Range("L37").Value = 30908.6
$ws.Range("M37").Value = -29527
$ws.Range("N37").Value = -31454.6
$ws.Range("H45").Value = 1173.4546
$ws.Range("I45").Value = 1003.75
$ws.Range("J45").Value = 1270.4286
$ws.Range("K45").Value = 1003.75
$ws.Range("L45").Value = 1270.4286
$ws.Range("M45").Value = -626.75
$ws.Range("N45").Value = -2024.4286
$ws.Range("H61").Value = 2150.0435
$ws.Range("I61").Value = 1912.6
$ws.Range("J61").Value = 3733
$ws.Range("K61").Value = 1912.6
$ws.Range("L61").Value = 3733
$ws.Range("M61").Value = -1700.6
$ws.Range("N61").Value = -4157
$ws.Range("H74").Value = 2030.2059
$ws.Range("I74").Value = 1518.8966
$ws.Range("J74").Value = 4995.8
$ws.Range("K74").Value = 1518.8966
$ws.Range("L74").Value = 4995.8
$ws.Range("M74").Value = -644.8966
$ws.Range("N74").Value = -6743.8
$ws.Range("H77").Value = 2030.2059
$ws.Range("I77").Value = 1518.8966
$ws.Range("J77").Value = 4995.8
$ws.Range("K77").Value = 7594.483
$ws.Range("L77").Value = 24979
$ws.Range("M77").Value = -3226.483
$ws.Range("N77").Value = -33715
$ws.Range("H116").Value = 644.1539
$ws.Range("I116").Value = 579.5454999999999
$ws.Range("J116").Value = 999.5
$ws.Range("K116").Value = 579.5454999999999
$ws.Range("L116").Value = 999.5
$ws.Range("M116").Value = 1714.4545
$ws.Range("N116").Value = -5587.5
$ws.Range("H136").Value = 2150.0435
$ws.Range("I136").Value = 1912.6
$ws.Range("J136").Value = 3733
$ws.Range("K136").Value = 5737.799999999999
$ws.Range("L136").Value = 11199
$ws.Range("M136").Value = -3187.799999999999
$ws.Range("N136").Value = -16299
$ws.Range("H137").Value = 40716.555
$ws.Range("J137").Value = 40716.555
$ws.Range("L137").Value = 40716.555
$ws.Range("N137").Value = -50916.555

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 644.1539
$ws.Range("I3").Value = 579.5454999999999
$ws.Range("J3").Value = 999.5
$ws.Range("K3").Value = 579.5454999999999
$ws.Range("L3").Value = 999.5
$ws.Range("M3").Value = -465.5454999999999
$ws.Range("N3").Value = -1227.5
$ws.Range("H134").Value = 2388.875
$ws.Range("I134").Value = 1411.7273
$ws.Range("J134").Value = 3924.3928
$ws.Range("K134").Value = 4235.1819
$ws.Range("L134").Value = 11773.1784
$ws.Range("M134").Value = -1700.1819
$ws.Range("N134").Value = -16843.1784
$ws.Range("H137").Value = 39546
$ws.Range("J137").Value = 40606.668
$ws.Range("L137").Value = 40606.668
$ws.Range("N137").Value = -50806.668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3900.875
$ws.Range("I132").Value = 2287.8572
$ws.Range("K132").Value = 6863.571599999999
$ws.Range("M132").Value = -4333.571599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 682.8946999999999
$ws.Range("I113").Value = 694.5
$ws.Range("J113").Value = 677.53845
$ws.Range("K113").Value = 2083.5
$ws.Range("L113").Value = 2032.61535
$ws.Range("M113").Value = 86.5
$ws.Range("N113").Value = -6372.61535
$ws.Range("H131").Value = 8475648
$ws.Range("J131").Value = 952.78845
$ws.Range("L131").Value = 2858.36535
$ws.Range("N131").Value = -12938.36535
$ws.Range("H137").Value = 2213.64
$ws.Range("I137").Value = 445
$ws.Range("J137").Value = 2772.158
$ws.Range("K137").Value = 1335
$ws.Range("L137").Value = 8316.474
$ws.Range("M137").Value = 3765
$ws.Range("N137").Value = -18516.474

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 31253738
$ws.Range("I80").Value = 83335000
$ws.Range("J80").Value = 4980
$ws.Range("K80").Value = 83335000
$ws.Range("L80").Value = 4980
$ws.Range("M80").Value = -83334002
$ws.Range("N80").Value = -6976
$ws.Range("H83").Value = 31253738
$ws.Range("I83").Value = 83335000
$ws.Range("J83").Value = 4980
$ws.Range("K83").Value = 416675000
$ws.Range("L83").Value = 24900
$ws.Range("M83").Value = -416670008
$ws.Range("N83").Value = -34884
$ws.Range("H107").Value = 6173425
$ws.Range("I107").Value = 368.8889
$ws.Range("J107").Value = 12346481
$ws.Range("K107").Value = 368.8889
$ws.Range("L107").Value = 12346481
$ws.Range("M107").Value = 1551.1111
$ws.Range("N107").Value = -12350321
$ws.Range("H113").Value = 1100
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1070
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 3410.8696
$ws.Range("I132").Value = 2136.8667
$ws.Range("J132").Value = 5799.625
$ws.Range("K132").Value = 6410.6001
$ws.Range("L132").Value = 17398.875
$ws.Range("M132").Value = -3880.6001
$ws.Range("N132").Value = -22458.875
$ws.Range("H137").Value = 61098.57
$ws.Range("J137").Value = 62106.152
$ws.Range("L137").Value = 62106.152
$ws.Range("N137").Value = -72306.152

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2009.6666
$ws.Range("J46").Value = 2700.2856
$ws.Range("L46").Value = 2700.2856
$ws.Range("N46").Value = -3076.2856
$ws.Range("H61").Value = 2266.6667
$ws.Range("I61").Value = 1900
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1900
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1698
$ws.Range("N61").Value = -3404
$ws.Range("H82").Value = 1562.1177
$ws.Range("I82").Value = 796.8570999999999
$ws.Range("J82").Value = 5133.3335
$ws.Range("K82").Value = 796.8570999999999
$ws.Range("L82").Value = 5133.3335
$ws.Range("M82").Value = -435.8570999999999
$ws.Range("N82").Value = -5855.3335
$ws.Range("H85").Value = 1562.1177
$ws.Range("I85").Value = 796.8570999999999
$ws.Range("J85").Value = 5133.3335
$ws.Range("K85").Value = 796.8570999999999
$ws.Range("L85").Value = 5133.3335
$ws.Range("M85").Value = 451.1429000000001
$ws.Range("N85").Value = -7629.3335
$ws.Range("H113").Value = 2266.6667
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1900
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 270
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 3827.3809
$ws.Range("I132").Value = 1718.6666
$ws.Range("K132").Value = 5155.9998
$ws.Range("M132").Value = -2625.9998
$ws.Range("H136").Value = 4175.35
$ws.Range("I136").Value = 2077.4614
$ws.Range("J136").Value = 8071.4287
$ws.Range("K136").Value = 6232.3842
$ws.Range("L136").Value = 24214.2861
$ws.Range("M136").Value = -3682.3842
$ws.Range("N136").Value = -29314.2861
$ws.Range("H140").Value = 67280.7
$ws.Range("J140").Value = 67280.7
$ws.Range("L140").Value = 67280.7
$ws.Range("N140").Value = -77640.7

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 125280000
$ws.Range("I62").Value = 250010000
$ws.Range("J62").Value = 550001.5
$ws.Range("K62").Value = 250010000
$ws.Range("L62").Value = 550001.5
$ws.Range("M62").Value = -250009376
$ws.Range("N62").Value = -551249.5
$ws.Range("H65").Value = 125280000
$ws.Range("I65").Value = 250010000
$ws.Range("J65").Value = 550001.5
$ws.Range("K65").Value = 1250050000
$ws.Range("L65").Value = 2750007.5
$ws.Range("M65").Value = -1250046880
$ws.Range("N65").Value = -2756247.5
$ws.Range("H107").Value = 611.5333000000001
$ws.Range("I107").Value = 446.4
$ws.Range("J107").Value = 1437.2
$ws.Range("K107").Value = 1339.2
$ws.Range("L107").Value = 4311.6
$ws.Range("M107").Value = 580.8000000000002
$ws.Range("N107").Value = -8151.6
$ws.Range("H126").Value = 2437.8438
$ws.Range("I126").Value = 1981.2354
$ws.Range("K126").Value = 5943.706200000001
$ws.Range("M126").Value = -3473.706200000001
$ws.Range("H128").Value = 42857.855
$ws.Range("J128").Value = 42857.855
$ws.Range("L128").Value = 42857.855
$ws.Range("N128").Value = -52817.855
$ws.Range("H131").Value = 50403.5
$ws.Range("J131").Value = 50403.5
$ws.Range("L131").Value = 50403.5
$ws.Range("N131").Value = -60483.5
$ws.Range("H138").Value = 46344.445
$ws.Range("J138").Value = 46344.445
$ws.Range("L138").Value = 46344.445
$ws.Range("N138").Value = -56624.445
